$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a "record" block for rows 3..11:
# D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion), S (Precio $/Kg), T (Kg/unidad)
$cols = @("D","L","M","N","O","P","Q","S","T")

# Snapshot current ("before") values for rows 3..11 first, since the update
# is a permutation of rows (not a simple overwrite) - read everything before
# writing anything.
$snapshot = @{}
for ($r = 3; $r -le 11; $r++) {
    $row = @{}
    foreach ($col in $cols) {
        $row[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $row
}

# Target row -> source row mapping (row's new data comes from the old data
# that used to live at the source row).
$mapping = @{
    3  = 10
    4  = 8
    5  = 11
    6  = 4
    7  = 9
    8  = 5
    9  = 3
    10 = 6
    11 = 7
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $srcData = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $srcData[$col]
    }
}
